$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E16:E27) previously listed periods in
# chronological order (2007, 2008, ..., 2012, 2101, ..., 2106).
# Old statement-of-account periods are removed and new ones are added,
# so the list is now shown in reverse order (2106, 2105, ..., 2007).
$ws.Range("E16").Value = "2106"
$ws.Range("E17").Value = "2105"
$ws.Range("E18").Value = "2104"
$ws.Range("E19").Value = "2103"
$ws.Range("E20").Value = "2102"
$ws.Range("E21").Value = "2101"
$ws.Range("E22").Value = "2012"
$ws.Range("E23").Value = "2011"
$ws.Range("E24").Value = "2010"
$ws.Range("E25").Value = "2009"
$ws.Range("E26").Value = "2008"
$ws.Range("E27").Value = "2007"
